$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.297706961631775
$ws.Range("B1").Value = 2.120767831802368
$ws.Range("C1").Value = 4.784872055053711
$ws.Range("D1").Value = 3.374036312103271
$ws.Range("E1").Value = 1.354151844978333
